$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6429
$ws.Range("E2").Value = 180
$ws.Range("F2").Value = 180
$ws.Range("G2").Value = 125
$ws.Range("H2").Value = 104
$ws.Range("I2").Value = 104
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 4658
$ws.Range("L2").Value = 2221
$ws.Range("M2").Value = 2437
$ws.Range("N2").Value = 2416
$ws.Range("O2").Value = 21
$ws.Range("P2").Value = 173
$ws.Range("Q2").Value = 227
$ws.Range("R2").Value = -374
$ws.Range("S2").Value = 75
$ws.Range("T2").Value = 40
$ws.Range("U2").Value = 187
$ws.Range("V2").Value = 1000
$ws.Range("W2").Value = 2.79
$ws.Range("X2").Value = 1.62
$ws.Range("Y2").Value = 4.71
$ws.Range("Z2").Value = 2.25
$ws.Range("AA2").Value = 91.15000000000001
$ws.Range("AB2").Value = 1351.13
$ws.Range("AC2").Value = 254
$ws.Range("AD2").Value = 24.14
$ws.Range("AE2").Value = 5286
$ws.Range("AF2").Value = 1.16
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 46704354

# Row 3
$ws.Range("D3").Value = 6505
$ws.Range("E3").Value = 234
$ws.Range("F3").Value = 234
$ws.Range("G3").Value = 185
$ws.Range("H3").Value = 134
$ws.Range("I3").Value = 132
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 4942
$ws.Range("L3").Value = 2350
$ws.Range("M3").Value = 2593
$ws.Range("N3").Value = 2569
$ws.Range("O3").Value = 23
$ws.Range("P3").Value = 173
$ws.Range("Q3").Value = -180
$ws.Range("R3").Value = 26
$ws.Range("S3").Value = 57
$ws.Range("T3").Value = 129
$ws.Range("U3").Value = -309
$ws.Range("V3").Value = 1100
$ws.Range("W3").Value = 3.59
$ws.Range("X3").Value = 2.06
$ws.Range("Y3").Value = 5.28
$ws.Range("Z3").Value = 2.79
$ws.Range("AA3").Value = 90.63
$ws.Range("AB3").Value = 1435.8
$ws.Range("AC3").Value = 282
$ws.Range("AD3").Value = 29.2
$ws.Range("AE3").Value = 5622
$ws.Range("AF3").Value = 1.46
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 46704354

# Row 4
$ws.Range("D4").Value = 6240
$ws.Range("E4").Value = 377
$ws.Range("F4").Value = 377
$ws.Range("G4").Value = 312
$ws.Range("H4").Value = 237
$ws.Range("I4").Value = 234
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 5087
$ws.Range("L4").Value = 2276
$ws.Range("M4").Value = 2811
$ws.Range("N4").Value = 2785
$ws.Range("O4").Value = 26
$ws.Range("P4").Value = 173
$ws.Range("Q4").Value = 656
$ws.Range("R4").Value = -114
$ws.Range("S4").Value = -36
$ws.Range("T4").Value = 20
$ws.Range("U4").Value = 636
$ws.Range("V4").Value = 1100
$ws.Range("W4").Value = 6.04
$ws.Range("X4").Value = 3.8
$ws.Range("Y4").Value = 8.74
$ws.Range("Z4").Value = 4.72
$ws.Range("AA4").Value = 80.95999999999999
$ws.Range("AB4").Value = 1555.38
$ws.Range("AC4").Value = 501
$ws.Range("AD4").Value = 15
$ws.Range("AE4").Value = 6093
$ws.Range("AF4").Value = 1.23
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 46704354

# Row 5
$ws.Range("D5").Value = 6243
$ws.Range("E5").Value = 342
$ws.Range("F5").Value = 342
$ws.Range("G5").Value = 339
$ws.Range("H5").Value = 249
$ws.Range("I5").Value = 249
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5340
$ws.Range("L5").Value = 2310
$ws.Range("M5").Value = 3030
$ws.Range("N5").Value = 3003
$ws.Range("O5").Value = 27
$ws.Range("P5").Value = 173
$ws.Range("Q5").Value = 599
$ws.Range("R5").Value = -102
$ws.Range("S5").Value = -270
$ws.Range("T5").Value = 28
$ws.Range("U5").Value = 570
$ws.Range("V5").Value = 900
$ws.Range("W5").Value = 5.48
$ws.Range("X5").Value = 3.99
$ws.Range("Y5").Value = 8.609999999999999
$ws.Range("Z5").Value = 4.78
$ws.Range("AA5").Value = 76.23999999999999
$ws.Range("AB5").Value = 1708.45
$ws.Range("AC5").Value = 534
$ws.Range("AD5").Value = 10.72
$ws.Range("AE5").Value = 6658
$ws.Range("AF5").Value = 0.86
$ws.Range("AG5").Value = 122
$ws.Range("AH5").Value = 2.13
$ws.Range("AI5").Value = 21.73
$ws.Range("AJ5").Value = 46704354

# Row 6
$ws.Range("D6").Value = 6429
$ws.Range("E6").Value = 340
$ws.Range("F6").Value = 340
$ws.Range("G6").Value = 307
$ws.Range("H6").Value = 223
$ws.Range("I6").Value = 222
$ws.Range("K6").Value = 6533
$ws.Range("L6").Value = 3343
$ws.Range("M6").Value = 3190
$ws.Range("N6").Value = 3162
$ws.Range("P6").Value = 173
$ws.Range("Q6").Value = 383
$ws.Range("R6").Value = -720
$ws.Range("S6").Value = 866
$ws.Range("T6").Value = 16
$ws.Range("U6").Value = 367
$ws.Range("V6").Value = 1850
$ws.Range("W6").Value = 5.29
$ws.Range("X6").Value = 3.47
$ws.Range("Y6").Value = 7.2
$ws.Range("Z6").Value = 3.76
$ws.Range("AA6").Value = 104.79
$ws.Range("AB6").Value = 1794.61
$ws.Range("AC6").Value = 475
$ws.Range("AD6").Value = 5.67
$ws.Range("AE6").Value = 7010
$ws.Range("AF6").Value = 0.38
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 46704354

# Row 7
$ws.Range("D7").Value = 6530
$ws.Range("E7").Value = 270
$ws.Range("G7").Value = 620
$ws.Range("H7").Value = 530
$ws.Range("I7").Value = 520
$ws.Range("K7").Value = 24910
$ws.Range("L7").Value = 18110
$ws.Range("M7").Value = 6800
$ws.Range("N7").Value = 6770
$ws.Range("P7").Value = 670
$ws.Range("Q7").Value = 170
$ws.Range("R7").Value = -17520
$ws.Range("S7").Value = 17020
$ws.Range("T7").Value = 70
$ws.Range("U7").ClearContents()
$ws.Range("W7").Value = 4.13
$ws.Range("X7").Value = 8.119999999999999
$ws.Range("Y7").Value = 10.47
$ws.Range("Z7").Value = 3.37
$ws.Range("AA7").Value = 266.32
$ws.Range("AC7").Value = 462
$ws.Range("AD7").Value = 4.83
$ws.Range("AE7").Value = 5107
$ws.Range("AF7").Value = 0.44
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 6870
$ws.Range("E8").Value = 330
$ws.Range("G8").Value = 310
$ws.Range("H8").Value = 230
$ws.Range("I8").Value = 230
$ws.Range("K8").Value = 9500
$ws.Range("L8").Value = 2470
$ws.Range("M8").Value = 7030
$ws.Range("N8").Value = 7000
$ws.Range("P8").Value = 670
$ws.Range("Q8").Value = 100
$ws.Range("R8").Value = 19310
$ws.Range("S8").Value = -15810
$ws.Range("T8").Value = 240
$ws.Range("U8").ClearContents()
$ws.Range("W8").Value = 4.8
$ws.Range("X8").Value = 3.35
$ws.Range("Y8").Value = 3.34
$ws.Range("Z8").Value = 1.34
$ws.Range("AA8").Value = 35.14
$ws.Range("AC8").Value = 171
$ws.Range("AD8").Value = 13.01
$ws.Range("AE8").Value = 5280
$ws.Range("AF8").Value = 0.42
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 7210
$ws.Range("E9").Value = 360
$ws.Range("G9").Value = 440
$ws.Range("H9").Value = 330
$ws.Range("I9").Value = 330
$ws.Range("K9").Value = 9850
$ws.Range("L9").Value = 2480
$ws.Range("M9").Value = 7370
$ws.Range("N9").Value = 7330
$ws.Range("P9").Value = 670
$ws.Range("Q9").Value = 80
$ws.Range("R9").Value = 220
$ws.Range("S9").Value = -50
$ws.Range("T9").Value = 240
$ws.Range("U9").ClearContents()
$ws.Range("W9").Value = 4.99
$ws.Range("X9").Value = 4.58
$ws.Range("Y9").Value = 4.61
$ws.Range("Z9").Value = 3.41
$ws.Range("AA9").Value = 33.65
$ws.Range("AC9").Value = 246
$ws.Range("AD9").Value = 9.07
$ws.Range("AE9").Value = 5529
$ws.Range("AF9").Value = 0.4
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AI9").ClearContents()
